# Added columns for path parameters
#
# Tests sheet: add "param:type" / "param:uuid" columns (I, J) and a new
# "Missing Required Param" validation test row.
# Documentation sheet: add a "Parameter Descriptions:" section describing
# the new path parameters, and a "Required parameters" note.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Tests"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New columns I and J, width ~20 (matches the other "wide" columns already
# present, just a bit wider since these hold a uuid/description text).
$ws1.Columns.Item(9).ColumnWidth = 19.14
$ws1.Columns.Item(10).ColumnWidth = 19.14

# New header cells
$ws1.Cells.Item(1, 9).Value = "param:type"
$ws1.Cells.Item(1, 10).Value = "param:uuid"

# Existing "Basic Test" row (row 2) gets values for the new columns
$ws1.Cells.Item(2, 9).Value = "work"
$ws1.Cells.Item(2, 10).Value = "12345678-1234-1234-1234-123456789abc"

# New row 3: a validation test exercising a missing required parameter
$ws1.Cells.Item(3, 1).Value = "delete-data - Missing Required Param"
$ws1.Cells.Item(3, 2).Value = "Test DELETE /data/:type/:uuid with missing required parameters"
# Leading apostrophe forces literal text (matches how "true"/"false" are
# stored as text elsewhere in this sheet rather than as booleans).
$ws1.Cells.Item(3, 3).Value = "'true"
$ws1.Cells.Item(3, 4).Value = 400
$ws1.Cells.Item(3, 5).Value = 10000
$ws1.Cells.Item(3, 6).Value = 2000
$ws1.Cells.Item(3, 7).Value = 500
$ws1.Cells.Item(3, 8).Value = "delete-data,validation"
$ws1.Cells.Item(3, 9).Value = ""
$ws1.Cells.Item(3, 10).Value = "12345678-1234-1234-1234-123456789abc"

# Extend the "number stored as text" error-suppression to cover the grown
# table (best effort - mirrors Range.Errors(xlNumberAsText).Ignore = True).
try {
    $ws1.Range("A1:J3").Errors.Item(9).Ignore = $true
} catch {}

# ---------------------------------------------------------------------
# Sheet 2: "Documentation"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Make room: 3 blank rows right after the (now repurposed) row 17, pushing
# "Endpoint-Specific Notes:" and friends further down the sheet.
$ws2.Range("A18:A20").EntireRow.Insert()

# One more blank row just above the final "Document operations..." note, to
# hold the new "Required parameters" line.
$ws2.Range("A25").EntireRow.Insert()

$ws2.Cells.Item(17, 1).Value = "Parameter Descriptions:"

$ws2.Cells.Item(18, 1).Value = "param:type"
$ws2.Cells.Item(18, 2).Value = "type parameter (string) (REQUIRED - highlighted in yellow)"

$ws2.Cells.Item(19, 1).Value = "param:uuid"
$ws2.Cells.Item(19, 2).Value = "uuid parameter (string) (REQUIRED - highlighted in yellow)"

# Row 20 stays blank (separator before "Endpoint-Specific Notes:").

$ws2.Cells.Item(25, 1).Value = "• Required parameters: type, uuid"

# Best effort, mirrors the Tests sheet.
try {
    $ws2.Range("A1:B26").Errors.Item(9).Ignore = $true
} catch {}
